# Apply cell-level edits per the commit diff.
# Rows 17-34 have their field values permuted among fixed row positions
# (records were re-ordered upstream); some optional columns are
# added/removed per row as a side effect.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range("A17").Value = 111359889
$ws.Range("B17").Value = 95538
$ws.Range("E17").Value = 221941
$ws.Range("F17").Value = "Plattlummer"
$ws.Range("G17").Value = "Lycopodium complanatum"
$ws.Range("H17").Value = "L."
$ws.Range("Q17").Value = 491909.0202035823
$ws.Range("R17").Value = 6785498.341940038

# Row 18
$ws.Range("A18").Value = 111356612
$ws.Range("B18").Value = 90666
$ws.Range("E18").Value = 4364
$ws.Range("F18").Value = "Dropptaggsvamp"
$ws.Range("G18").Value = "Hydnellum ferrugineum"
$ws.Range("H18").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q18").Value = 491951.0498785287
$ws.Range("R18").Value = 6785511.741186595

# Row 19
$ws.Range("A19").Value = 111491649
$ws.Range("AC19").Value = "mycel"
$ws.Range("B19").Value = 90666
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 4364
$ws.Range("F19").Value = "Dropptaggsvamp"
$ws.Range("G19").Value = "Hydnellum ferrugineum"
$ws.Range("H19").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q19").Value = 491979.6153062462
$ws.Range("R19").Value = 6785548.307010972

# Row 20
$ws.Range("A20").Value = 111491680
$ws.Range("B20").Value = 56414
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 100049
$ws.Range("F20").Value = "Spillkråka"
$ws.Range("G20").Value = "Dryocopus martius"
$ws.Range("H20").Value = "(Linnaeus, 1758)"
$ws.Range("M20").Value = "äldre spår"

# Row 21
$ws.Range("A21").Value = 111491681
$ws.Range("AC21").Value = "äldre fruktkreopp"
$ws.Range("Q21").Value = 491929.8523854768
$ws.Range("R21").Value = 6785530.587422797

# Row 22
$ws.Range("A22").Value = 111491641
$ws.Range("AC22").Value = "Minst 4 platser i området"
$ws.Range("B22").Value = 90666
$ws.Range("E22").Value = 4364
$ws.Range("F22").Value = "Dropptaggsvamp"
$ws.Range("G22").Value = "Hydnellum ferrugineum"
$ws.Range("H22").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I22").Value = "'4"

# Row 23
$ws.Range("A23").Value = 111491685
$ws.Range("B23").Value = 88819
$ws.Range("E23").Value = 5685
$ws.Range("F23").Value = "Gullgröppa"
$ws.Range("G23").Value = "Pseudomerulius aureus"
$ws.Range("H23").Value = "(Fr.) Jülich"
$ws.Range("Q23").Value = 491909.4940688942
$ws.Range("R23").Value = 6785494.484901348
$ws.Range("S23").Value = 5

# Row 25
$ws.Range("A25").Value = 111491657
$ws.Range("AC25").Value = "vid bohål"
$ws.Range("B25").Value = 89425
$ws.Range("E25").Value = 5442
$ws.Range("F25").Value = "Tallticka"
$ws.Range("G25").Value = "Porodaedalea pini"
$ws.Range("H25").Value = "(Brot.) Murrill"
$ws.Range("Q25").Value = 491946.35724353
$ws.Range("R25").Value = 6785570.554389503

# Row 26
$ws.Range("A26").Value = 111491639
$ws.Range("AC26").Value = "Spridd och riklig i området"
$ws.Range("B26").Value = 8377
$ws.Range("E26").Value = 106545
$ws.Range("F26").Value = "Mindre märgborre"
$ws.Range("G26").Value = "Tomicus minor"
$ws.Range("H26").Value = "(Hartig, 1834)"
$ws.Range("M26").Value = "äldre gnagspår"
$ws.Range("Q26").Value = 491993.9996831641
$ws.Range("R26").Value = 6785505.377163783
$ws.Range("S26").Value = 100

# Row 27
$ws.Range("A27").Value = 111612720
$ws.Range("B27").Value = 88924
$ws.Range("D27").Value = "LC"
$ws.Range("E27").Value = 256703
$ws.Range("F27").Value = "Tallfingersvamp"
$ws.Range("G27").Value = "Ramaria eosanguinea"
$ws.Range("H27").Value = "R.H.Petersen"
$ws.Range("I27").Value = "'2"
$ws.Range("Q27").Value = 491993.9996831641
$ws.Range("R27").Value = 6785505.377163783
$ws.Range("S27").Value = 100

# Row 28
$ws.Range("A28").Value = 111612736
$ws.Range("B28").Value = 56398
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 100109
$ws.Range("F28").Value = "Tretåig hackspett"
$ws.Range("G28").Value = "Picoides tridactylus"
$ws.Range("H28").Value = "(Linnaeus, 1758)"
$ws.Range("M28").Value = "äldre spår"
$ws.Range("Q28").Value = 491952.3910193561
$ws.Range("R28").Value = 6785464.984647369
$ws.Range("S28").Value = 10

# Row 29
$ws.Range("A29").Value = 111612738
$ws.Range("B29").Value = 56414
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"

# Row 32
$ws.Range("A32").Value = 111682658
$ws.Range("B32").Value = 90689
$ws.Range("E32").Value = 5966
$ws.Range("F32").Value = "Motaggsvamp"
$ws.Range("G32").Value = "Sarcodon squamosus"
$ws.Range("H32").Value = "(Schaeff.) Quél."

# Row 33
$ws.Range("A33").Value = 111682665
$ws.Range("AC33").Value = "Längs stigen/traktorspåret strax utanför gränsmarkeringen som syns på träden"
$ws.Range("B33").Value = 90682
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 2059
$ws.Range("F33").Value = "Skrovlig taggsvamp"
$ws.Range("G33").Value = "Hydnellum scabrosum"
$ws.Range("H33").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q33").Value = 492024.0709204427
$ws.Range("R33").Value = 6785567.485207787

# Row 34
$ws.Range("A34").Value = 111682769
$ws.Range("AJ34").Value = "gran"
$ws.Range("AK34").Value = "Picea abies"
$ws.Range("AO34").Value = "Picea abies"
$ws.Range("B34").Value = 89980
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 1179
$ws.Range("F34").Value = "Gräddticka"
$ws.Range("G34").Value = "Perenniporia subacida"
$ws.Range("H34").Value = "(Peck) Donk"
$ws.Range("Q34").Value = 491952.3910193561
$ws.Range("R34").Value = 6785464.984647369

# Cells that no longer hold any value in the target rows
$ws.Range("L18").ClearContents()
$ws.Range("AC20").ClearContents()
$ws.Range("AF20").ClearContents()
$ws.Range("J20").ClearContents()
$ws.Range("L22").ClearContents()
$ws.Range("M22").ClearContents()
$ws.Range("AC23").ClearContents()
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("L27").ClearContents()
$ws.Range("M27").ClearContents()
$ws.Range("AF28").ClearContents()
$ws.Range("J28").ClearContents()
$ws.Range("AJ33").ClearContents()
$ws.Range("AK33").ClearContents()
$ws.Range("AO33").ClearContents()
$ws.Range("AC34").ClearContents()

# Cells that become present-but-empty placeholders; clear any existing
# value first, then copy an empty placeholder cell from the same row
# (column K is never otherwise touched) so the target cell exists
# without holding a real value.
$ws.Range("L17").ClearContents()
$ws.Range("K17").Copy($ws.Range("L17"))
$ws.Range("L20").ClearContents()
$ws.Range("K20").Copy($ws.Range("L20"))
$ws.Range("I23").ClearContents()
$ws.Range("K23").Copy($ws.Range("I23"))
$ws.Range("AF25").ClearContents()
$ws.Range("K25").Copy($ws.Range("AF25"))
$ws.Range("J25").ClearContents()
$ws.Range("K25").Copy($ws.Range("J25"))
$ws.Range("L26").ClearContents()
$ws.Range("K26").Copy($ws.Range("L26"))
$ws.Range("AF27").ClearContents()
$ws.Range("K27").Copy($ws.Range("AF27"))
$ws.Range("J27").ClearContents()
$ws.Range("K27").Copy($ws.Range("J27"))
$ws.Range("I28").ClearContents()
$ws.Range("K28").Copy($ws.Range("I28"))
$ws.Range("L28").ClearContents()
$ws.Range("K28").Copy($ws.Range("L28"))
